$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2117.6667
$ws.Range("I32").Value = 2726.5
$ws.Range("K32").Value = 2726.5
$ws.Range("M32").Value = -2400.5
$ws.Range("H34").Value = 27666.666
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 27666.666
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 27666.666
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -28072.666
$ws.Range("H36").Value = 27666.666
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 27666.666
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 27666.666
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -29096.666
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 817.8182
$ws.Range("I132").Value = 879.6
$ws.Range("K132").Value = 2638.8
$ws.Range("M132").Value = -108.8000000000002
$ws.Range("H138").Value = 2482.611
$ws.Range("I138").Value = 1421.2858
$ws.Range("K138").Value = 4263.857400000001
$ws.Range("M138").Value = 876.1425999999992

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 546
$ws.Range("I41").Value = 546
$ws.Range("K41").Value = 546
$ws.Range("M41").Value = -132
$ws.Range("H43").Value = 37999
$ws.Range("I43").Value = 37999
$ws.Range("K43").Value = 37999
$ws.Range("M43").Value = -37686
$ws.Range("H102").Value = 16232640
$ws.Range("I102").Value = 1223835.2
$ws.Range("K102").Value = 1223835.2
$ws.Range("M102").Value = -1222213.2
$ws.Range("H122").Value = 1966.6666
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -10900
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 5477.75
$ws.Range("I75").Value = 5477.75
$ws.Range("K75").Value = 5477.75
$ws.Range("M75").Value = -4541.75
$ws.Range("H78").Value = 5477.75
$ws.Range("I78").Value = 5477.75
$ws.Range("K78").Value = 16433.25
$ws.Range("M78").Value = -11753.25
$ws.Range("H80").Value = 1328.6666
$ws.Range("I80").Value = 1393
$ws.Range("K80").Value = 1393
$ws.Range("M80").Value = -395
$ws.Range("H83").Value = 1328.6666
$ws.Range("I83").Value = 1393
$ws.Range("K83").Value = 6965
$ws.Range("M83").Value = -1973
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744
$ws.Range("H107").Value = 83069.8
$ws.Range("I107").Value = 102837.25
$ws.Range("K107").Value = 102837.25
$ws.Range("M107").Value = -100917.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5715399.5
$ws.Range("I6").Value = 10000200
$ws.Range("J6").Value = 2333.3333
$ws.Range("K6").Value = 10000200
$ws.Range("L6").Value = 2333.3333
$ws.Range("M6").Value = -10000087
$ws.Range("N6").Value = -2559.3333
$ws.Range("H31").Value = 2757.875
$ws.Range("I31").Value = 1325.8422
$ws.Range("J31").Value = 8199.6
$ws.Range("K31").Value = 1325.8422
$ws.Range("L31").Value = 8199.6
$ws.Range("M31").Value = -1030.8422
$ws.Range("N31").Value = -8789.6
$ws.Range("H34").Value = 2757.875
$ws.Range("I34").Value = 1325.8422
$ws.Range("J34").Value = 8199.6
$ws.Range("K34").Value = 1325.8422
$ws.Range("L34").Value = 8199.6
$ws.Range("M34").Value = -1123.8422
$ws.Range("N34").Value = -8603.6
$ws.Range("H86").Value = 100006200
$ws.Range("I86").Value = 142858860
$ws.Range("K86").Value = 142858860
$ws.Range("M86").Value = -142857737
$ws.Range("H89").Value = 100006200
$ws.Range("I89").Value = 142858860
$ws.Range("K89").Value = 714294300
$ws.Range("M89").Value = -714288684
$ws.Range("H132").Value = 1998.2
$ws.Range("I132").Value = 2097.75
$ws.Range("K132").Value = 6293.25
$ws.Range("M132").Value = -3763.25
$ws.Range("H134").Value = 999.6667
$ws.Range("I134").Value = 999.6667
$ws.Range("K134").Value = 2999.0001
$ws.Range("M134").Value = -464.0001000000002

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 499.5
$ws.Range("I8").Value = 499.5
$ws.Range("K8").Value = 1498.5
$ws.Range("M8").Value = -1359.5
$ws.Range("H34").Value = 2554.2104
$ws.Range("I34").Value = 425.5
$ws.Range("J34").Value = 2631.6182
$ws.Range("K34").Value = 1276.5
$ws.Range("L34").Value = 7894.8546
$ws.Range("M34").Value = -1192.5
$ws.Range("N34").Value = -8062.8546
$ws.Range("H39").Value = 1599.5333
$ws.Range("J39").Value = 1599.5333
$ws.Range("L39").Value = 4798.5999
$ws.Range("N39").Value = -5386.5999
$ws.Range("H55").Value = 3081.0952
$ws.Range("I55").Value = 1578
$ws.Range("J55").Value = 3682.3333
$ws.Range("K55").Value = 4734
$ws.Range("L55").Value = 11046.9999
$ws.Range("M55").Value = -4557
$ws.Range("N55").Value = -11400.9999
$ws.Range("H115").Value = 1031
$ws.Range("J115").Value = 1031
$ws.Range("L115").Value = 3093
$ws.Range("N115").Value = -5443
$ws.Range("H140").Value = 2924
$ws.Range("I140").Value = 2913.1428
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 8739.428400000001
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -3559.428400000001
$ws.Range("N140").Value = -19360

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4745.706
$ws.Range("I122").Value = 3562.7273
$ws.Range("K122").Value = 10688.1819
$ws.Range("M122").Value = -8238.1819

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5224
$ws.Range("H82").Value = 2951.4375
$ws.Range("I82").Value = 1804.4445
$ws.Range("K82").Value = 1804.4445
$ws.Range("M82").Value = -1443.4445
$ws.Range("H85").Value = 2951.4375
$ws.Range("I85").Value = 1804.4445
$ws.Range("K85").Value = 1804.4445
$ws.Range("M85").Value = -556.4445000000001
$ws.Range("H122").Value = 3482.4119
$ws.Range("I122").Value = 3378.6428
$ws.Range("K122").Value = 10135.9284
$ws.Range("M122").Value = -7685.928400000001
$ws.Range("H132").Value = 9399.799999999999
$ws.Range("I132").Value = 9500
$ws.Range("K132").Value = 28500
$ws.Range("M132").Value = -25970
$ws.Range("H136").Value = 2225
$ws.Range("J136").Value = 2633.3333
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("N136").Value = -12999.9999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4743
$ws.Range("I126").Value = 3822.375
$ws.Range("K126").Value = 11467.125
$ws.Range("M126").Value = -8997.125
